$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two input parameters that drive the whole recalculation
$ws.Range("F1").Value = 1.0009999999999999
$ws.Range("H1").Value = 20

# Add new descriptive header labels (new shared strings) next to the formulas
$ws.Range("F2").Value = "1/RTT1"
$ws.Range("G2").Value = "1/RTT2"
$ws.Range("H2").Value = "bottleneck rate"

# Move the active selection
$ws.Range("I8").Select() | Out-Null
